$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4 and 3 entirely (delete bottom-up so row numbers stay valid)
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# --- Row 2 updates ---

# A2: numeric-looking id string "02250002" -> "02250003".
# Using a formula that evaluates to text, then collapsing it to a static
# value via copy/paste-values keeps the cell a text (shared-string) cell
# instead of Excel auto-coercing a pure-digit string into a number.
$ws.Range("A2").Formula = "=""02250003"""
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)

# B2: date/time serial value update
$ws.Range("B2").Value = 45712.41756395833

# E2: equipment label update
$ws.Range("E2").Value = "ECLAIRAGE EXTERIEUR"

# I2: description now populated
$ws.Range("I2").Value = "Mise en maintenance"

# K2: supervisor id resolved to a generic label
$ws.Range("K2").Value = "user id"

# L2 (Date de cloture) and N2 (Edite par) are cleared out entirely
# (Clear, not ClearContents, so the cell -- and its style -- is fully
# removed instead of leaving an empty but still-present/styled cell)
$ws.Range("L2").Clear()
$ws.Range("N2").Clear()

# M2: creator id resolved to a different id
$ws.Range("M2").Value = "0237a803-e675-49df-9d7b-25f2b329704b"

# O2: status update
$ws.Range("O2").Value = "EN ATTENTE"

$excel.CutCopyMode = 0
